$wb = $excel.ActiveWorkbook

# --- 1. Append the new Proximity event row (row 13) -----------------------
$proximity = $wb.Worksheets.Item("Proximity")

# Force column A to stay plain text ("2026-01-28" would otherwise be
# auto-converted to a date serial number by Excel's type inference).
$proximity.Range("A13").NumberFormat = "@"
$proximity.Range("A13").Value = "2026-01-28"
$proximity.Range("B13").Value = "17:49:07"
$proximity.Range("C13").Value = "17:00"
$proximity.Range("D13").Value = "Bathroom Door"
$proximity.Range("E13").Value = "ENTER"
$proximity.Range("F13").Value = "User ENTERED Bathroom"

# --- 2. Add the new "Sleep" sheet, fill it in, then relocate it after
#        "Camera" (the sheet handle is repositioned last, since moving it
#        earlier leaves the variable pointing at the wrong sheet) ----------
$sleep = $wb.Worksheets.Add()
$sleep.Name = "Sleep"

# Header row
$sleep.Range("A1").Value = "Date"
$sleep.Range("B1").Value = "Timestamp"
$sleep.Range("C1").Value = "Hour"
$sleep.Range("D1").Value = "Location"
$sleep.Range("E1").Value = "BedState"
$sleep.Range("F1").Value = "HeartRate"
$sleep.Range("G1").Value = "BreathRate"
$sleep.Range("H1").Value = "Status"

# Keep the Date column as text for every data row too.
$sleep.Range("A2:A9").NumberFormat = "@"

$sleepData = @(
    @("2026-01-28", "17:49:00", "17:00", "Bedroom", "In Bed", 0,   0,  "Occupied"),
    @("2026-01-28", "17:49:00", "17:00", "Bedroom", "In Bed", 108, 60, "Occupied"),
    @("2026-01-28", "17:49:02", "17:00", "Bedroom", "In Bed", 92,  44, "Occupied"),
    @("2026-01-28", "17:49:02", "17:00", "Bedroom", "In Bed", 56,  8,  "Occupied"),
    @("2026-01-28", "17:49:04", "17:00", "Bedroom", "In Bed", 50,  2,  "Occupied"),
    @("2026-01-28", "17:49:05", "17:00", "Bedroom", "In Bed", 51,  3,  "Occupied"),
    @("2026-01-28", "17:49:07", "17:00", "Bedroom", "In Bed", 50,  2,  "Occupied"),
    @("2026-01-28", "17:49:16", "17:00", "Bedroom", "In Bed", 49,  1,  "Occupied")
)

$r = 2
foreach ($row in $sleepData) {
    $sleep.Cells.Item($r, 1).Value = $row[0]
    $sleep.Cells.Item($r, 2).Value = $row[1]
    $sleep.Cells.Item($r, 3).Value = $row[2]
    $sleep.Cells.Item($r, 4).Value = $row[3]
    $sleep.Cells.Item($r, 5).Value = $row[4]
    $sleep.Cells.Item($r, 6).Value = $row[5]
    $sleep.Cells.Item($r, 7).Value = $row[6]
    $sleep.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Now relocate the finished sheet to sit right after "Camera".
$camera = $wb.Worksheets.Item("Camera")
$sleep.Move($null, $camera)
